$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update Version value (row 3, column B): 0.1.1 -> 0.2.0
$ws.Cells.Item(3, 2).Value = "0.2.0"

# 2. Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2023-10-19T17:05:12+00:00"

# 3. Insert a new row after "Contact" (row 10) for "Jurisdiction" / "iso:code:3166:FR",
#    pushing "Description" and all following rows down by one.
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row that is now below the newly inserted row
# (the former row 11, "Description", now at row 12) onto the new row so the
# new row keeps the same style (borders / wrap / vertical alignment) as the
# rest of the data rows.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"
